$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "cheval_numero" row (row 6) from its old position in the
#    "cheval" group. This shifts rows 7-11 up to rows 6-10, preserving
#    their existing formatting (styles move with the rows).
$ws.Rows(6).Delete()

# 2. Re-add "cheval_numero" as a new standalone row at the bottom (row 11),
#    matching the style used by the "pari" rows (A6:E8) rather than its
#    former "cheval" group style.
$ws.Range("A8:E8").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A11").Value = "cheval_numero"

# 3. Add the new "cheval_numero" column (F), mirroring column D's header
#    style and each row's data style.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "cheval_numero"

$ws.Range("D2:D11").Copy()
$ws.Range("F2:F11").PasteSpecial(-4122)

# Size the new column to fit its header text, same as the existing columns.
$ws.Columns("F").AutoFit()

# 4. Update the selection to mirror the author's saved selection.
$ws.Range("F17").Select()
